# DATE_plate_layout.xlsx -- 2020-03-03 growth plate experiment layout
#
# Replaces the generic placeholder grid on the "strain" sheet with the
# actual dilution-series plate layout (5 strain/condition rows x 10
# dilution steps + a "blank" row), resizes the used range down to
# A1:K6, sets the column widths used by that sheet, and restores the
# original cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strain")

# ---- new grid contents -----------------------------------------------
$row1 = @("none_HG105_0.0", "none_HG105_0.1", "none_HG105_0.2", "none_HG105_0.3", "none_HG105_0.4", "none_HG105_0.5", "none_HG105_0.6", "none_HG105_0.7", "none_HG105_0.8", "none_HG105_0.9", "none_HG105_1.0")
$row2 = @("O2_R1740_0.0", "O2_R1740_0.1", "O2_R1740_0.2", "O2_R1740_0.3", "O2_R1740_0.4", "O2_R1740_0.5", "O2_R1740_0.6", "O2_R1740_0.7", "O2_R1740_0.8", "O2_R1740_0.9", "O2_R1740_1.0")
$row3 = @("O2_R260_0.0", "O2_R260_0.1", "O2_R260_0.2", "O2_R260_0.3", "O2_R260_0.4", "O2_R260_0.5", "O2_R260_0.6", "O2_R260_0.7", "O2_R260_0.8", "O2_R260_0.9", "O2_R260_1.0")
$row4 = @("O2_R22_0.0", "O2_R22_0.1", "O2_R22_0.2", "O2_R22_0.3", "O2_R22_0.4", "O2_R22_0.5", "O2_R22_0.6", "O2_R22_0.7", "O2_R22_0.8", "O2_R22_0.9", "O2_R22_1.0")
$row5 = @("O2_R0_0.0", "O2_R0_0.1", "O2_R0_0.2", "O2_R0_0.3", "O2_R0_0.4", "O2_R0_0.5", "O2_R0_0.6", "O2_R0_0.7", "O2_R0_0.8", "O2_R0_0.9", "O2_R0_1.0")
$row6 = @("blank", "blank", "blank", "blank", "blank", "blank", "blank", "blank", "blank", "blank", "blank")

$rows = @($row1, $row2, $row3, $row4, $row5, $row6)

for ($r = 1; $r -le 6; $r++) {
    $vals = $rows[$r - 1]
    for ($c = 1; $c -le 11; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# ---- drop the now-unused tail of the old 12x8 placeholder grid -------
$ws.Range("A7:L8").ClearContents()
$ws.Range("L1:L6").ClearContents()

# ---- column widths matching the new, narrower layout ------------------
$ws.Columns.Item(1).ColumnWidth = 24.5
$ws.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws.Columns.Item(3).ColumnWidth = 21.666666666666668
$ws.Columns.Item(4).ColumnWidth = 25.166666666666668
$ws.Columns.Item(5).ColumnWidth = 20.666666666666668
$ws.Columns.Item(6).ColumnWidth = 20.666666666666668
$ws.Columns.Item(7).ColumnWidth = 14.666666666666666
$ws.Columns.Item(8).ColumnWidth = 16.166666666666668
$ws.Columns.Item(9).ColumnWidth = 16.833333333333332
$ws.Columns.Item(10).ColumnWidth = 40.666666666666664
$ws.Columns.Item(11).ColumnWidth = 37.333333333333336

# ---- restore the active selection on the sheet -------------------------
$ws.Range("B44").Select() | Out-Null
